$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H41").Value = 1284.9
$ws.Range("J41").Value = 2492
$ws.Range("L41").Value = 2492
$ws.Range("N41").Value = -3372
$ws.Range("H106").Value = 0
$ws.Range("I106").Value = 0
$ws.Range("K106").Value = 0
$ws.Range("M106").ClearContents()
$ws.Range("H132").Value = 2871799.5
$ws.Range("I132").Value = 2991424.5
$ws.Range("K132").Value = 8974273.5
$ws.Range("M132").Value = -8971743.5
$ws.Range("H137").Value = 11321.093
$ws.Range("I137").Value = 18182.904
$ws.Range("K137").Value = 54548.712
$ws.Range("M137").Value = -51998.712
$ws.Range("H138").Value = 9196.154
$ws.Range("J138").Value = 9850.25
$ws.Range("L138").Value = 29550.75
$ws.Range("N138").Value = -39830.75

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1006.1667
$ws.Range("I2").Value = 964.7143
$ws.Range("J2").Value = 1064.2
$ws.Range("K2").Value = 964.7143
$ws.Range("L2").Value = 1064.2
$ws.Range("M2").Value = -851.7143
$ws.Range("N2").Value = -1290.2
$ws.Range("H35").Value = 8045.6665
$ws.Range("I35").Value = 4974.8
$ws.Range("J35").Value = 23400
$ws.Range("K35").Value = 4974.8
$ws.Range("L35").Value = 23400
$ws.Range("M35").Value = -4568.8
$ws.Range("N35").Value = -24212
$ws.Range("H74").Value = 178619.17
$ws.Range("I74").Value = 208743.17
$ws.Range("K74").Value = 208743.17
$ws.Range("M74").Value = -207869.17
$ws.Range("H77").Value = 178619.17
$ws.Range("I77").Value = 208743.17
$ws.Range("K77").Value = 1043715.85
$ws.Range("M77").Value = -1039347.85
$ws.Range("H116").Value = 1006.1667
$ws.Range("I116").Value = 964.7143
$ws.Range("J116").Value = 1064.2
$ws.Range("K116").Value = 964.7143
$ws.Range("L116").Value = 1064.2
$ws.Range("M116").Value = 1329.2857
$ws.Range("N116").Value = -5652.2

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1006.1667
$ws.Range("I3").Value = 964.7143
$ws.Range("J3").Value = 1064.2
$ws.Range("K3").Value = 964.7143
$ws.Range("L3").Value = 1064.2
$ws.Range("M3").Value = -850.7143
$ws.Range("N3").Value = -1292.2
$ws.Range("H19").Value = 0
$ws.Range("J19").Value = 0
$ws.Range("L19").Value = 0
$ws.Range("N19").ClearContents()
$ws.Range("H20").Value = 40434.46
$ws.Range("I20").Value = 51961.8
$ws.Range("J20").Value = 2010
$ws.Range("K20").Value = 51961.8
$ws.Range("L20").Value = 2010
$ws.Range("M20").Value = -51714.8
$ws.Range("N20").Value = -2504
$ws.Range("H88").Value = 65166.668
$ws.Range("J88").Value = 65166.668
$ws.Range("L88").Value = 65166.668
$ws.Range("N88").Value = -65978.66800000001
$ws.Range("H91").Value = 65166.668
$ws.Range("J91").Value = 65166.668
$ws.Range("L91").Value = 65166.668
$ws.Range("N91").Value = -67974.66800000001
$ws.Range("H97").Value = 14942.556
$ws.Range("I97").Value = 2474.5
$ws.Range("J97").Value = 24917
$ws.Range("K97").Value = 2474.5
$ws.Range("L97").Value = 24917
$ws.Range("M97").Value = -1483.5
$ws.Range("N97").Value = -26899

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5595.844
$ws.Range("I31").Value = 3731.48
$ws.Range("J31").Value = 7926.3
$ws.Range("K31").Value = 3731.48
$ws.Range("L31").Value = 7926.3
$ws.Range("M31").Value = -3436.48
$ws.Range("N31").Value = -8516.299999999999
$ws.Range("H34").Value = 5595.844
$ws.Range("I34").Value = 3731.48
$ws.Range("J34").Value = 7926.3
$ws.Range("K34").Value = 3731.48
$ws.Range("L34").Value = 7926.3
$ws.Range("M34").Value = -3529.48
$ws.Range("N34").Value = -8330.299999999999
$ws.Range("H68").Value = 31250
$ws.Range("H71").Value = 31250
$ws.Range("H134").Value = 3166.7026
$ws.Range("I134").Value = 2866.742
$ws.Range("J134").Value = 4716.5
$ws.Range("K134").Value = 8600.226000000001
$ws.Range("L134").Value = 14149.5
$ws.Range("M134").Value = -6065.226000000001
$ws.Range("N134").Value = -19219.5

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 2281.5
$ws.Range("I5").Value = 1547.25
$ws.Range("K5").Value = 4641.75
$ws.Range("M5").Value = -4529.75
$ws.Range("H39").Value = 4372.5
$ws.Range("J39").Value = 7500
$ws.Range("L39").Value = 22500
$ws.Range("N39").Value = -23088
$ws.Range("H135").Value = 2281.5
$ws.Range("I135").Value = 1547.25
$ws.Range("K135").Value = 13925.25
$ws.Range("M135").Value = -11390.25

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 340.27777
$ws.Range("J2").Value = 810
$ws.Range("L2").Value = 810
$ws.Range("N2").Value = -1036
$ws.Range("H97").Value = 648.71875
$ws.Range("J97").Value = 356.2857
$ws.Range("L97").Value = 356.2857
$ws.Range("N97").Value = -1348.2857

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H13").Value = 10438
$ws.Range("J13").Value = 1500
$ws.Range("L13").Value = 1500
$ws.Range("N13").Value = -1780
$ws.Range("H22").Value = 3656.6667
$ws.Range("I22").Value = 2583.8333
$ws.Range("K22").Value = 2583.8333
$ws.Range("M22").Value = -2288.8333
$ws.Range("H27").Value = 3656.6667
$ws.Range("I27").Value = 2583.8333
$ws.Range("K27").Value = 2583.8333
$ws.Range("M27").Value = -2476.8333
$ws.Range("H61").Value = 1907.3077
$ws.Range("I61").Value = 1890.7273
$ws.Range("K61").Value = 1890.7273
$ws.Range("M61").Value = -1688.7273
$ws.Range("H98").Value = 60000
$ws.Range("J98").Value = 60000
$ws.Range("L98").Value = 60000
$ws.Range("N98").Value = -65990
$ws.Range("H113").Value = 1907.3077
$ws.Range("I113").Value = 1890.7273
$ws.Range("K113").Value = 1890.7273
$ws.Range("M113").Value = 279.2727
$ws.Range("H136").Value = 3638.6875
$ws.Range("I136").Value = 2609.7778
$ws.Range("K136").Value = 7829.3334
$ws.Range("M136").Value = -5279.3334

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 1075.6428
$ws.Range("J107").Value = 1283
$ws.Range("L107").Value = 3849
$ws.Range("N107").Value = -7689
$ws.Range("H123").Value = 143333.33
$ws.Range("J123").Value = 154000
$ws.Range("L123").Value = 154000
$ws.Range("N123").Value = -163800
$ws.Range("H126").Value = 169593.17
$ws.Range("I126").Value = 2324.875
$ws.Range("K126").Value = 6974.625
$ws.Range("M126").Value = -4504.625
